$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell content, text must be inserted in this exact order so that
# --- the resulting shared-string table indices line up with the target file:
# --- B15 (BEGINER) -> 30, B16 (ADVANCED) -> 31, B14 (Header) -> 32,
# --- C14 (bullets) -> 33, C15 (myheader example) -> 34, C16 (insertion_sort example) -> 35

$b15 = @'
## Create my header (BEGINER)
'@

$b16 = @'
## Create my header (ADVANCED)
'@

$b14 = @'
## Header
'@

$c14 = @'
* #include <file> //compile look for header in __SYSTEM__ DIRECTORY
* #include "file"  //coompiler look for header in __CURRENT__ DIRECOTRY
* It is ok to compile without header impl, as long as the impl of header is present when program is called
* The compilation is __stupid__, merely just copy the headers into the target program...
'@

$c15 = @'
### Prepare myheader.h and its content is:
```c
int add(int a,int b)
{
return(a+b);
}
```
### Prepare the main program main.c:
```c
#include<stdio.h>
#include"myhead.h"
void main() {
   int num1 = 10, num2 = 10, num3;
   num3 = add(num1, num2);
   printf("Addition of Two numbers : %d", num3);
}
```
### Compile 
`cc main.c`  (must ensure myheader.h at the current folder)
### Run
`./a.exe`
'@

$c16 = @'
### Prepare header file "insertion_sort.h"
```c
//insertion_sort.h file
#ifndef _insertion_sort_h
#define _insertion_sort_h
/*Sorts an integer array. Takes a pointer to the first element and the length of the array as input. 
Returns 0 on successful sort.*/
int insertion_sort(int*, int); 
#endif
```
### Prepare header impl "insertion_sort.c"
```c
//insertion_sort.c file
#include "insertion_sort.h"
int insertion_sort(int *a, int n){
 int i=1;
 int j=i;
 int t;
 for(i=1; i<n; ++i){
  for(j=0; j<i; ++j){
   if(a[i]<a[j]){
    t = a[i];
    a[i] = a[j];
    a[j] = t;
   }
  }
 }
 return 0;
}
```
### Prepare the main "insertion_main.c"
```c
//insertion_main.c file
#include <stdio.h>
#include "insertion_sort.h"
void print_array(int*, int);
int main(){
 int a[] = {10, 9, 8, 7, 6, 5, 4, 3, 2, 1};
 insertion_sort(a, 10);
 print_array(a, 10);
}
void print_array(int *arr, int n){
 int i;
 for(i=0; i<n; ++i){
  printf("%d ", arr[i]);
 }
}
```
### Compile 
```bash
cc -c insertion_sort.c #this will give you insertion_sort.o
cc -c insertion_main.c #this will give you insertion_main.o
cc -o insertion_main insertion_main.o insertion_sort.o #this will create the insertion_main executable
```
### Run
`./insertion_main`
'@

# Column A repeats the existing "# C" language label (shared string already present).
$ws.Range("A14").Value = "# C"
$ws.Range("A15").Value = "# C"
$ws.Range("A16").Value = "# C"

# Insert in the order that reproduces the target shared-string ordering.
$ws.Range("B15").Value = $b15
$ws.Range("B16").Value = $b16
$ws.Range("B14").Value = $b14
$ws.Range("C14").Value = $c14
$ws.Range("C15").Value = $c15
$ws.Range("C16").Value = $c16

# Match the row heights used throughout the sheet.
$ws.Rows.Item(14).RowHeight = 47.25
$ws.Rows.Item(15).RowHeight = 47.25
$ws.Rows.Item(16).RowHeight = 47.25

# Match the view state recorded in the target file (active cell / selection;
# this headless COM host does not persist pane scroll position / topLeftCell).
$ws.Range("B24").Select()
